# Natmi following Dr Hou advice
# Rebuild the Inhbb-Acvr1 LR-pair table with the "ECs" sending cluster added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2-7) completely; we'll rewrite rows 2-10 below.
$ws.Range("A2:T7").Clear()

$data = @(
    @("ECs",  "Inhbb", "Acvr1", "ECs",  2, 0.6666666666666666, 3.345805333333333, 10.037416,       0.3489465220682754, 0.3489465220682754, 3, 1, 4.695610666666666,  14.086832,  0.1802066564018305, 0.1802066564018305, 15.71059921179022, 141.395392906112,  0.06288248600497147, 0.06288248600497146),
    @("ECs",  "Inhbb", "Acvr1", "FAPs", 2, 0.6666666666666666, 3.345805333333333, 10.037416,       0.3489465220682754, 0.3489465220682754, 3, 1, 15.51448033333333,  46.543441,  0.5954098039960916, 0.5954098039960916, 51.90843104316177, 467.1758793884559, 0.2077661803097897,  0.2077661803097897),
    @("ECs",  "Inhbb", "Acvr1", "sCs",  2, 0.6666666666666666, 3.345805333333333, 10.037416,       0.3489465220682754, 0.3489465220682754, 3, 1, 5.846719333333333,  17.540158,  0.2243835396020779, 0.2243835396020779, 19.56198472796977, 176.057862551728,  0.07829785575351422, 0.07829785575351421),
    @("FAPs", "Inhbb", "Acvr1", "ECs",  3, 1,                  5.480061666666667, 16.440185,       0.5715360783999618, 0.5715360783999618, 3, 1, 4.695610666666666,  14.086832,  0.1802066564018305, 0.1802066564018305, 25.73223601599111, 231.59012414392,   0.1029946057014716,  0.1029946057014716),
    @("FAPs", "Inhbb", "Acvr1", "FAPs", 3, 1,                  5.480061666666667, 16.440185,       0.5715360783999618, 0.5715360783999618, 3, 1, 15.51448033333333,  46.543441,  0.5954098039960916, 0.5954098039960916, 85.0203089529539,  765.1827805765851, 0.3402981844168161,  0.3402981844168161),
    @("FAPs", "Inhbb", "Acvr1", "sCs",  3, 1,                  5.480061666666667, 16.440185,       0.5715360783999618, 0.5715360783999618, 3, 1, 5.846719333333333,  17.540158,  0.2243835396020779, 0.2243835396020779, 32.04038249435889, 288.36344244923,   0.1282432882816741,  0.1282432882816741),
    @("sCs",  "Inhbb", "Acvr1", "ECs",  3, 1,                  0.7624369999999999, 2.287311,       0.07951739953176286,0.07951739953176286, 3, 1, 4.695610666666666,  14.086832,  0.1802066564018305, 0.1802066564018305, 3.580107309861333, 32.220965788752,   0.01432956469538747, 0.01432956469538747),
    @("sCs",  "Inhbb", "Acvr1", "FAPs", 3, 1,                  0.7624369999999999, 2.287311,       0.07951739953176286,0.07951739953176286, 3, 1, 15.51448033333333,  46.543441,  0.5954098039960916, 0.5954098039960916, 11.82881384190567, 106.459324577151,  0.04734543926948583, 0.04734543926948583),
    @("sCs",  "Inhbb", "Acvr1", "sCs",  3, 1,                  0.7624369999999999, 2.287311,       0.07951739953176286,0.07951739953176286, 3, 1, 5.846719333333333,  17.540158,  0.2243835396020779, 0.2243835396020779, 4.457755148348665, 40.11979633513799, 0.01784239556688956, 0.01784239556688956)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $ws.Cells.Item($row, 12).Value = $rec[11]
    $ws.Cells.Item($row, 13).Value = $rec[12]
    $ws.Cells.Item($row, 14).Value = $rec[13]
    $ws.Cells.Item($row, 15).Value = $rec[14]
    $ws.Cells.Item($row, 16).Value = $rec[15]
    $ws.Cells.Item($row, 17).Value = $rec[16]
    $ws.Cells.Item($row, 18).Value = $rec[17]
    $ws.Cells.Item($row, 19).Value = $rec[18]
    $ws.Cells.Item($row, 20).Value = $rec[19]
    $row = $row + 1
}
